$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 72, shifting the existing rows 72-104 down to 73-105.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data entry.
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44845
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = 100112031
$ws.Range("G72").Value = "Poroto verde"
$ws.Range("H72").Value = "Magnum"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 35
$ws.Range("K72").Value = 35000
$ws.Range("L72").Value = 35000
$ws.Range("M72").Value = 35000
$ws.Range("N72").Value = "`$/malla 25 kilos"
$ws.Range("O72").Value = "Perú"
$ws.Range("P72").Value = 1400
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"
